# "Schedule Management" workbook correction:
# The WBS numbering for the last task group (previously mis-labelled as
# Sprint "4.x" while already living after another "Sprint 4" block) is
# renumbered to Sprint "5.x" — i.e. rows 33-38, column A go from
# 4 / 4.1 / 4.2 / 4.3 / 4.4 / 4.5  ->  5 / 5.1 / 5.2 / 5.3 / 5.4 / 5.5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 5
$ws.Range("A34").Value = 5.0999999999999996
$ws.Range("A35").Value = 5.2
$ws.Range("A36").Value = 5.3
$ws.Range("A37").Value = 5.4
$ws.Range("A38").Value = 5.5

# Reflect where the user ended up after making the edit: scrolled further
# down the sheet, with cell A38 selected (was topLeftCell A7 / selection
# C12 before the edit).
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A38").Select()
